# Applies the cryptos list update described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Donor cell with the default (unstyled) cell format used by all data rows,
# so that forcing text via NumberFormat can be reverted without leaving a stray style index.
$defaultStyleCell = $ws.Range("B2")

function Set-TextValue {
    param($range, [string]$text)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = $defaultStyleCell.Style
}

Set-TextValue $ws.Range("D2") "26.133.18"
Set-TextValue $ws.Range("D3") "1.654.55"
$ws.Range("E3").Value = "  -0.11%  "
Set-TextValue $ws.Range("D4") "1.002"
$ws.Range("E4").Value = "  -0.22%  "
Set-TextValue $ws.Range("D5") "218.61"
$ws.Range("E5").Value = "  +0.03%  "
Set-TextValue $ws.Range("D6") "0.5238"
$ws.Range("E6").Value = "  -0.14%  "
$ws.Range("E7").Value = "  -0.23%  "
$ws.Range("E8").Value = "  +1.31%  "
$ws.Range("E9").Value = "  +1.04%  "
$ws.Range("E10").Value = "  -0.07%  "
Set-TextValue $ws.Range("D11") "0.07699"
$ws.Range("E11").Value = "  -1.48%  "
Set-TextValue $ws.Range("D12") "4.631"
$ws.Range("E12").Value = "  +3.32%  "
Set-TextValue $ws.Range("D13") "1.652.15"
$ws.Range("E13").Value = "  -0.61%  "
Set-TextValue $ws.Range("D14") "1.882.53"
$ws.Range("E14").Value = "  -0.14%  "
Set-TextValue $ws.Range("D15") "0.5614"
$ws.Range("E15").Value = "  +1.20%  "
Set-TextValue $ws.Range("D16") "0.0₅8190"
$ws.Range("E16").Value = "  +2.10%  "
Set-TextValue $ws.Range("D17") "65.46"
$ws.Range("E17").Value = "  +0.77%  "
Set-TextValue $ws.Range("D18") "26.130.92"
$ws.Range("E18").Value = "  -0.06%  "
$ws.Range("E19").Value = "  -0.17%  "
Set-TextValue $ws.Range("D20") "4.652"
$ws.Range("E20").Value = "  +0.51%  "
Set-TextValue $ws.Range("D21") "10.48"
$ws.Range("E21").Value = "  +4.06%  "
Set-TextValue $ws.Range("D22") "191.99"
$ws.Range("E22").Value = "  -1.38%  "
Set-TextValue $ws.Range("D23") "5.953"
$ws.Range("E23").Value = "  -0.02%  "
Set-TextValue $ws.Range("D24") "1.003"
$ws.Range("E24").Value = "  -0.21%  "
Set-TextValue $ws.Range("D25") "145.03"
$ws.Range("E25").Value = "  -1.31%  "
Set-TextValue $ws.Range("D26") "0.1194"
$ws.Range("E26").Value = "  -0.97%  "
$ws.Range("E27").Value = "  +1.36%  "
Set-TextValue $ws.Range("D28") "15.93"
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("E29").Value = "  +1.61%  "
Set-TextValue $ws.Range("D30") "0.05444"
$ws.Range("E30").Value = "  -4.46%  "
Set-TextValue $ws.Range("D31") "1.271"
$ws.Range("E31").Value = "  +0.07%  "
$ws.Range("E32").Value = "  -0.59%  "
Set-TextValue $ws.Range("D33") "3.370"
$ws.Range("E33").Value = "  +0.86%  "
Set-TextValue $ws.Range("D34") "1.559"
$ws.Range("E34").Value = "  -1.79%  "
Set-TextValue $ws.Range("D35") "0.9521"
$ws.Range("E35").Value = "  +0.24%  "
$ws.Range("E36").Value = "  -0.87%  "
Set-TextValue $ws.Range("D37") "2.401"
$ws.Range("E37").Value = "  -0.60%  "
Set-TextValue $ws.Range("D38") "0.5667"
$ws.Range("E38").Value = "  -0.61%  "
Set-TextValue $ws.Range("D39") "0.01581"
$ws.Range("E39").Value = "  -0.87%  "
$ws.Range("E40").Value = "  -1.24%  "
$ws.Range("E41").Value = "  -0.20%  "
Set-TextValue $ws.Range("D42") "0.8355"
$ws.Range("E42").Value = "  -1.20%  "
Set-TextValue $ws.Range("D43") "1.028.71"
$ws.Range("E43").Value = "  -3.33%  "
Set-TextValue $ws.Range("D44") "101.21"
$ws.Range("E44").Value = "  -2.15%  "
Set-TextValue $ws.Range("D45") "1.793.43"
$ws.Range("E45").Value = "  -0.12%  "
Set-TextValue $ws.Range("D46") "57.75"
$ws.Range("E46").Value = "  +0.07%  "
Set-TextValue $ws.Range("D47") "0.9998"
$ws.Range("E47").Value = "  -0.67%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue $ws.Range("D48") "0.0₈103"
$ws.Range("E48").Value = "  -0.10%  "
Set-TextValue $ws.Range("D49") "7.999"
$ws.Range("E49").Value = "  +0.23%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws.Range("D50") "0.4338"
$ws.Range("E50").Value = "  -1.44%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D51") "0.05191"
$ws.Range("E51").Value = "  -2.20%  "

Write-Host "Applied 94 cell updates"
